# Cập nhật data.xlsx từ công cụ QR
# A new QR scan result is inserted as the newest record at the top of the
# data table (row 2), pushing the previous newest record down to row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (existing row 2 shifts down to row 3).
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the newly scanned record.
$ws.Cells.Item(2, 1).Value2  = "30em912l8wj"
$ws.Cells.Item(2, 2).Value2  = "pretbgsw"
$ws.Cells.Item(2, 3).Value2  = "Hộ kinh doanh"
$ws.Cells.Item(2, 4).Value2  = "Madam Thu Bakery, 21C, Võ Văn Tần, Ninh Kiều, Ninh Kiều District, Cần Thơ, 94111, Vietnam"
$ws.Cells.Item(2, 5).Value2  = "https://www.google.com/maps/search/?api=1&query=10.032100,105.786400"
$ws.Cells.Item(2, 6).Value2  = "2025-08-21T07:30:22.697Z"
$ws.Cells.Item(2, 7).Value2  = ""
$ws.Cells.Item(2, 8).Value2  = ""
$ws.Cells.Item(2, 9).Value2  = ""
$ws.Cells.Item(2, 10).Value2 = ""
$ws.Cells.Item(2, 11).Value2 = ""
$ws.Cells.Item(2, 12).Value2 = "Nguyễn Văn B"
$ws.Cells.Item(2, 13).Value2 = "3243bcc63be98bf8"
$ws.Cells.Item(2, 14).Value2 = "1441c319972f62a62bfa8570c3f4f56c988ddea5dcec0830037775df5e4fcf85"
